$d = $word.ActiveDocument

# Locate the "2. Shareholders" paragraph (last paragraph of the original body)
$lastIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs($lastIndex)
$anchorRange = $anchor.Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

# Build each logical block of the comments section as its own paragraph first,
# then splice the paragraph marks back out so everything lands in one <w:p> as
# separate runs (matching how Word keeps each block as a distinct run).
$paraIndexes = @()

# --- block 0 ---
$pIndex0 = $d.Paragraphs.Count
$p0 = $d.Paragraphs($pIndex0)
$r0 = $p0.Range
$r0.Collapse(0)
$insStart0 = $r0.Start
$r0.InsertAfter('=== Comments Section ===')
$insEnd0 = $r0.End
$r0.Collapse(0)
$r0.InsertParagraphAfter()
$paraIndexes += , @($insStart0, $insEnd0)

# --- block 1 ---
$pIndex1 = $d.Paragraphs.Count
$p1 = $d.Paragraphs($pIndex1)
$r1 = $p1.Range
$r1.Collapse(0)
$insStart1 = $r1.Start
$r1.InsertAfter([char]11 + 'The following issues were detected:' + [char]11)
$insEnd1 = $r1.End
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$paraIndexes += , @($insStart1, $insEnd1)

# --- block 2 ---
$pIndex2 = $d.Paragraphs.Count
$p2 = $d.Paragraphs($pIndex2)
$r2 = $p2.Range
$r2.Collapse(0)
$insStart2 = $r2.Start
$r2.InsertAfter([char]11 + '1. **Section:** Signatory Section' + [char]11 + '   - **Issue:** Missing signatory section or improper formatting.' + [char]11 + '   - **Severity:** High' + [char]11 + '   - **Suggestion:** Add a signed section with all required parties per ADGM template guidelines.' + [char]11)
$insEnd2 = $r2.End
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$paraIndexes += , @($insStart2, $insEnd2)

# --- block 3 ---
$pIndex3 = $d.Paragraphs.Count
$p3 = $d.Paragraphs($pIndex3)
$r3 = $p3.Range
$r3.Collapse(0)
$insStart3 = $r3.Start
$r3.InsertAfter([char]11 + '2. **Section:** General Clauses' + [char]11 + '   - **Issue:** Based on the provided ADGM regulations, I can identify a potential issue with the text.' + [char]11 + [char]11 + 'The text states that the jurisdiction is "ADGM courts." However, according to the regulations, the body corporate must have unlimited liability in its home jurisdiction. Since the text does not specify the home jurisdiction, it is unclear whether the members of the body corporate have unlimited liability in that jurisdiction.' + [char]11 + [char]11 + 'Additionally, the text does not specify the type of company being formed (e.g., unlimited company, limited company). According to the regulations, if the members have unlimited liability, the body corporate must apply for continuance as an unlimited company.' + [char]11 + [char]11 + 'Therefore, the text may be missing essential information, making it potentially invalid under ADGM regulations.' + [char]11 + '   - **Severity:** Medium' + [char]11 + '   - **Suggestion:** Review and align with ADGM templates.' + [char]11)
$insEnd3 = $r3.End
$paraIndexes += , @($insStart3, $insEnd3)

# Make the header block ("=== Comments Section ===") bold.
$boldRange = $d.Range($paraIndexes[0][0], $paraIndexes[0][1])
$boldRange.Bold = 1

# Merge the helper paragraphs back together by deleting the paragraph marks
# we introduced between each block (but not the very first boundary, which
# was already in the document). Walk backwards so earlier offsets remain valid.
$mergeEnd = $paraIndexes[2][1] + 1
$mergeRange = $d.Range($paraIndexes[2][1], $mergeEnd)
$mergeRange.Delete()
$mergeEnd = $paraIndexes[1][1] + 1
$mergeRange = $d.Range($paraIndexes[1][1], $mergeEnd)
$mergeRange.Delete()
$mergeEnd = $paraIndexes[0][1] + 1
$mergeRange = $d.Range($paraIndexes[0][1], $mergeEnd)
$mergeRange.Delete()

Write-Host "done"